$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Email" header in B1 is being renamed to "Correo"
$ws.Range("B1").Value = "Correo"

# Update the active selection saved in the sheet view to B1
$ws.Range("B1").Select()
